$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# EFT block (rows 3-6) — fill empty GARANTI (C) and HALKBANK (H) cells
$ws.Range("C3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("C4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("C5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

$ws.Range("C6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("H6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

# HAVALE block (rows 8-11)
$ws.Range("C8").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("H8").Value = "15,23 TL - 30,47 TL - 304,72 TL"

$ws.Range("C9").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("H9").Value = "15,23 TL - 30,47 TL - 304,72 TL"

$ws.Range("C10").Value = "14,29 TL - 28,57 TL - 300 TL"
$ws.Range("H10").Value = "15,23 TL - 30,47 TL - 304,72 TL"

$ws.Range("C11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("H11").Value = "3,05 TL - 6,1 TL - 76,18 TL"

# SWIFT block (rows 12-14)
$ws.Range("C12").Value = "WU: 1.000,01 USD–9,51 USD"

$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 909,5 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 1.114 TL"
$ws.Range("H13").Value = "Hesaba: Asgari 1 TL | Azami 6,09 TL"

$ws.Range("C14").Value = "40.000 TL - 1.904,76 TL"
$ws.Range("H14").Value = "2.100 TL - 4.300 TL"
